# Administration schema update:
# - USERS table: add "salt" (row9) and "admin" (row15) columns, shifting
#   fname/lname/comp/phone/reg_date down by one row; rename "company" -> "comp"
# - CHANGES table: add "change_date" column (row9)
# - MESSAGE_ATTACH (renamed from MESSAGES_ATTACH) table: reorder FK columns
#   (tid/uid instead of uid/mid-header mixup fix), add message, create_date_m,
#   image, icon, create_date rows
# - TEXTS table: "h"/"t" -> true/false

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- USERS table: insert "salt" row after password (row8), before fname ---
$ws.Range("C9").Value = "salt"
$ws.Range("D9").Value = "NN"

$ws.Range("C10").Value = "fname"
$ws.Range("D10").Value = "NN"

$ws.Range("C11").Value = "lname"
$ws.Range("D11").Value = "NN"

$ws.Range("C12").Value = "comp"
$ws.Range("D12").Value = "NN"

$ws.Range("C13").Value = "phone"
$ws.Range("D13").ClearContents()

$ws.Range("C14").Value = "reg_date"
$ws.Range("D14").Value = "CURRENT TIME"

$ws.Range("C15").Value = "admin"
$ws.Range("D15").Value = "NN"

# --- CHANGES table: add change_date column ---
$ws.Range("S9").Value = "change_date"
$ws.Range("T9").Value = "CURRENT TIME"

# --- TEXTS table: page flag values rename ---
$ws.Range("D24").Value = "true/false"

# --- MESSAGE_ATTACH header rename (was MESSAGES_ATTACH) ---
$ws.Range("S22").Value = "MESSAGE_ATTACH"

# --- TOPICS table rows 24-25 swap (uid/header order fix) ---
$ws.Range("K24").Value = "uid"
$ws.Range("L24").Value = "FK NN"
$ws.Range("K25").Value = "header"
$ws.Range("L25").Value = "NN"

# --- MESSAGES table rows 24-27 ---
$ws.Range("O24").Value = "tid"
$ws.Range("P24").Value = "FK NN"
$ws.Range("O25").Value = "uid"
$ws.Range("P25").Value = "FK NN"
$ws.Range("O26").Value = "message"
$ws.Range("P26").Value = "NN"
$ws.Range("O27").Value = "create_date_m"
$ws.Range("P27").Value = "CURRENT TIME"

# --- MESSAGE_ATTACH table rows 28-30 (image/icon/create_date) ---
$ws.Range("C28").Value = "image"
$ws.Range("D28").Value = "NN"
$ws.Range("C29").Value = "icon"
$ws.Range("C30").Value = "create_date"
$ws.Range("D30").Value = "CURRENT TIME"

# --- Column T width (newly used range extends to column T) ---
# ColumnWidth uses Excel's character-width metric; 13.1666... yields the same
# stored OOXML width (14 "screen" units) as the sheet's other bestFit columns.
$ws.Columns.Item(20).ColumnWidth = 13.1666666667

# --- Update dimension / selection to match new used range ---
$ws.Range("S23").Select()
